$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet originally held 2 data rows (rows 2-3). This edit:
#   - replaces row 2 with a new incident (SR_PTO016M_HNI / Phuc Tho)
#   - replaces row 3 with what used to be row 2's data (UL_TTT093M_HNI)
#   - appends two brand-new rows (4-5) for the STY003 station exclusion
#     (commit message: "loai tru tram sty003")
# Build rows 4 & 5 by copying row 3's full formatting (borders, fill,
# alignment, wrap) before updating any values, so the new rows inherit
# the same style as the rest of the table in one shot.
# ------------------------------------------------------------------

$ws.Range("A3:AC3").Copy($ws.Range("A4:AC4"))
$ws.Range("A3:AC3").Copy($ws.Range("A5:AC5"))

# Row 2 - new incident record
$ws.Range("F2").Value = "SR_PTO016M_HNI"
$ws.Range("G2").Value = "Tam-Hiep-Thon-Thuong-PTO_HNI"
$ws.Range("J2").Value = "POWER_AC_EAS"
$ws.Range("L2").Value = "09/05/2025 14:07:23"
$ws.Range("T2").Value = "Phúc Thọ"
$ws.Range("AA2").Value = "Trạm viễn thông loại 2"

# Row 3 - now holds the record previously shown in row 2
$ws.Range("F3").Value = "UL_TTT093M_HNI"
$ws.Range("G3").Value = "THACH-HOA-TTT_HNI"
$ws.Range("J3").Value = "POWER_AC_EAS"
$ws.Range("L3").Value = "09/05/2025 12:39:56"
$ws.Range("T3").Value = "Thạch Thất"
$ws.Range("V3").Value = ""
$ws.Range("AA3").Value = "Trạm viễn thông loại 3"

# Row 4 - new STY003 (4G) record
$ws.Range("F4").Value = "4G-STY003M-HNI"
$ws.Range("G4").Value = "Lang-Van-Hoa-STY_HNI"
$ws.Range("J4").Value = "SITE_OOS"
$ws.Range("L4").Value = "08/05/2025 08:56:47"
$ws.Range("T4").Value = "Sơn Tây"
$ws.Range("V4").Value = "184602- VTHN ĐKTĐ - Thay cột treo anten trạm , dự kiến từ 08h00 ngày 07/05 đến ngày 12/05 - 4 - hanhhh - 08/05/2025 09:17:21"
$ws.Range("AA4").Value = "Trạm viễn thông loại 1"

# Row 5 - new STY003 (3G) record
$ws.Range("F5").Value = "3G_STY003M_HNI"
$ws.Range("G5").Value = "Lang-Van-Hoa-STY_HNI"
$ws.Range("J5").Value = "SITE_OOS"
$ws.Range("L5").Value = "08/05/2025 08:10:05"
$ws.Range("T5").Value = "Sơn Tây"
$ws.Range("V5").Value = "184602- VTHN ĐKTĐ - Thay cột treo anten trạm , dự kiến từ 08h00 ngày 07/05 đến ngày 12/05  - 1 - hanhhh - 08/05/2025 08:54:42"
$ws.Range("AA5").Value = "Trạm viễn thông loại 1"

# Column width tweaks: "Tên gợi nhớ" (G) narrower, "Tỉnh ghi chú" (V) much
# wider to fit the long maintenance notes now stored there.
$ws.Columns.Item(7).ColumnWidth = 29.8
$ws.Columns.Item(22).ColumnWidth = 126.8
